# Updated cryptos list with latest price/volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "43.188.76"
$ws.Cells.Item(2, 5).Value = "  -4.71%  "

$ws.Cells.Item(3, 4).Value = "2.236.72"
$ws.Cells.Item(3, 5).Value = "  -5.51%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "318.80"
$ws.Cells.Item(5, 5).Value = "  +1.64%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "101.12"
$ws.Cells.Item(6, 5).Value = "  -6.21%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.587"

$ws.Cells.Item(8, 5).Value = "  -0.16%  "

$ws.Cells.Item(9, 5).Value = "  -7.90%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "37.06"
$ws.Cells.Item(10, 5).Value = "  -9.07%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "54.49"
$ws.Cells.Item(11, 5).Value = "  -2.57%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0828"
$ws.Cells.Item(12, 5).Value = "  -9.67%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "7.69"
$ws.Cells.Item(13, 5).Value = "  -9.24%  "

$ws.Cells.Item(14, 5).Value = "  -1.33%  "

$ws.Cells.Item(15, 4).Value = "2.576.37"
$ws.Cells.Item(15, 5).Value = "  -5.77%  "

$ws.Cells.Item(16, 5).Value = "  -11.83%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "14.39"
$ws.Cells.Item(17, 5).Value = "  -6.02%  "

$ws.Cells.Item(18, 4).Value = "2.241.58"
$ws.Cells.Item(18, 5).Value = "  -5.82%  "

$ws.Cells.Item(19, 4).Value = "43.116.92"
$ws.Cells.Item(19, 5).Value = "  -4.89%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "14.38"
$ws.Cells.Item(20, 5).Value = "  -7.17%  "

$ws.Cells.Item(21, 5).Value = "  -8.90%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.53"
$ws.Cells.Item(22, 5).Value = "  -10.22%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "65.40"
$ws.Cells.Item(23, 5).Value = "  -10.78%  "

$ws.Cells.Item(24, 5).Value = "  -11.21%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "238.49"
$ws.Cells.Item(25, 5).Value = "  -8.59%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.16"
$ws.Cells.Item(26, 5).Value = "  -8.19%  "

$ws.Cells.Item(27, 5).Value = "  -0.25%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "10.07"
$ws.Cells.Item(29, 5).Value = "  -9.55%  "

$ws.Cells.Item(30, 5).Value = "  -2.42%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "6.40"
$ws.Cells.Item(31, 5).Value = "  -14.81%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "35.55"
$ws.Cells.Item(32, 5).Value = "  -3.89%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "20.48"
$ws.Cells.Item(33, 5).Value = "  -8.04%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.0878"
$ws.Cells.Item(34, 5).Value = "  -9.07%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "153.49"
$ws.Cells.Item(35, 5).Value = "  -7.87%  "

$ws.Cells.Item(36, 5).Value = "  -4.84%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "3.11"
$ws.Cells.Item(37, 5).Value = "  +5.95%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.96"
$ws.Cells.Item(38, 5).Value = "  +4.08%  "

$ws.Cells.Item(39, 5).Value = "  -6.86%  "

$ws.Cells.Item(40, 5).Value = "  -5.06%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.104"
$ws.Cells.Item(41, 5).Value = "  -10.90%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.70"
$ws.Cells.Item(42, 5).Value = "  -6.20%  "

$ws.Cells.Item(43, 5).Value = "  -8.09%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "12.90"
$ws.Cells.Item(44, 5).Value = "  -1.38%  "

$ws.Cells.Item(45, 5).Value = "  -0.16%  "

$ws.Cells.Item(46, 4).Value = "1.803.22"
$ws.Cells.Item(46, 5).Value = "  -0.66%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "87.41"
$ws.Cells.Item(47, 5).Value = "  -11.31%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.206"
$ws.Cells.Item(48, 5).Value = "  -9.40%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "76.74"
$ws.Cells.Item(49, 5).Value = "  -7.14%  "

$ws.Cells.Item(50, 5).Value = "  -10.19%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "59.28"
$ws.Cells.Item(51, 5).Value = "  -15.37%  "
